$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Item + Cost between row 2 (Chang) and row 4 (Aniseed Syrup),
# leaving Purchased quantities untouched.
$ws.Range("A2").Value = "Aniseed Syrup"
$ws.Range("C2").Value = 50

$ws.Range("A4").Value = "Chang"
$ws.Range("C4").Value = 95
